# The <id> tags for p100v_1 and p100v_2 were previously split across three
# separate runs ("<id>", the bare id value, "</id>") so that the id value
# could carry different (plainer) formatting than the surrounding tag
# markup. Collapse each of those triples back into a single run/text node
# now that the whole tag is meant to share one formatting run.

$d = $word.ActiveDocument

$targets = @("p100v_1", "p100v_2")

foreach ($id in $targets) {
    $old = "<id>" + $id + "</id>"
    $new = "<id>" + $id + "</id>"

    $find = $d.Content.Find
    $find.ClearFormatting()
    $ok = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Output "$id -> $ok"
}
